# standardize q by catchment area, now in mm/day
$wb = $excel.ActiveWorkbook

$wsWatersheds = $wb.Worksheets.Item("watersheds")
$wsHydrometric = $wb.Worksheets.Item("hydrometric")

# Watersheds tab: relabel the area columns and switch the conversion
# from Ha -> km^2 to Ha -> m^2 (C column stays hectares, D becomes m^2).
$wsWatersheds.Range("C1").Value = "Area (km^2)"
$wsWatersheds.Range("D1").Value = "Area (m^2)"

$wsWatersheds.Range("D2").Formula = "=C2*10000"

$wsWatersheds.Range("C3").Value = 1567.29
$wsWatersheds.Range("D3").Formula = "=C3*10000"

# Selection / active-sheet bookkeeping to match the authored state.
$wsHydrometric.Range("B2").Select()

$wsWatersheds.Activate()
$wsWatersheds.Range("D3").Select()
